$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SpawnList")

# Insert two new rows right after row 12 (after "c_dungeon"), pushing the
# remaining data down by 2 rows, then fill them in.
$ws.Rows.Item(13).Resize(2).Insert()

$ws.Cells.Item(13, 1).Value = "c_dungeon_forest"
$ws.Cells.Item(13, 2).Value = "EA 23.214"

$ws.Cells.Item(14, 1).Value = "c_machine"
$ws.Cells.Item(14, 2).Value = "EA 23.214"

# Append a new row after the last existing row ("darksoup", now at row 49
# after the insert above) with the new spawn list entry.
$ws.Cells.Item(50, 1).Value = "wreck_junk"
$ws.Cells.Item(50, 2).Value = "EA 23.228"
